# Apply the edit described by the diff:
#  - Set C1 to FALSE (it was TRUE)
#  - Delete rows 2-4 (Lampada da sala / Ar da sala / Ar do quarto) entirely
#  - Resulting used range/dimension becomes A1:C1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C1 value to FALSE
$ws.Range("C1").Value = $false

# Delete rows 2 through 4 (shifts cells up, removing the rows completely)
$ws.Range("A2:C4").EntireRow.Delete()
